$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new hires (e_cell / admission team) - "Dr. NAFEZA E" (row 9) and
# "Mr. CHAKKARAVARTHY KUMARESAN R" (row 10) - need to be moved up above the
# existing "Mrs. JENIFER HEPZIBA I" / "Mrs. JAIASHI J" rows (7-8), i.e. the
# two 2-row blocks (7:8) and (9:10) swap places, values+formatting included.
# Row 11 is untouched.

# Stage both blocks in a scratch area first so neither copy reads data that
# has already been overwritten by the other.
$ws.Range("A7:J8").Copy($ws.Range("A1000:J1001"))
$ws.Range("A9:J10").Copy($ws.Range("A1002:J1003"))

# Fully clear the live rows so that cells blank in the source don't leave
# stale values/styles behind when pasted over.
$ws.Range("A7:J10").Clear()

# Paste the blocks back in swapped order.
$ws.Range("A1002:J1003").Copy($ws.Range("A7:J8"))
$ws.Range("A1000:J1001").Copy($ws.Range("A9:J10"))

# Remove the scratch data entirely.
$ws.Range("A1000:J1003").Clear()

# Range.Copy doesn't carry per-row heights, so restore them explicitly.
$ws.Rows(7).RowHeight = 19.5
$ws.Rows(8).RowHeight = 72
$ws.Rows(9).RowHeight = 47.25
$ws.Rows(10).RowHeight = 47.25

# Match the saved view/selection state (user had scrolled to/selected row 8).
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("A8:XFD8").Select()
